$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object "object[,]" 29,18
$arr[0,0] = 251706
$arr[0,1] = "T3"
$arr[0,2] = 0
$arr[0,3] = 50.79365079365079
$arr[0,4] = "2025-05-12 07:00:00"
$arr[0,5] = "2025-05-12 07:00:00"
$arr[0,6] = "2025-05-12 07:00:00"
$arr[0,7] = "2025-05-12 07:50:47"
$arr[0,8] = 3200
$arr[0,9] = "foglio"
$arr[0,10] = "T3"
$arr[0,11] = 0
$arr[0,12] = 0
$arr[0,13] = "39764 (esterno)"
$arr[0,14] = "X"
$arr[0,15] = 39764
$arr[0,16] = "2025-05-14 00:00:00"
$arr[0,17] = 0
$arr[1,0] = 251455
$arr[1,1] = "BIMEC 2"
$arr[1,2] = 19
$arr[1,3] = 82.765625
$arr[1,4] = "2025-05-07 07:00:00"
$arr[1,5] = "2025-05-07 07:19:00"
$arr[1,6] = "2025-05-07 07:19:00"
$arr[1,7] = "2025-05-07 08:41:45"
$arr[1,8] = 5297
$arr[1,9] = "bobina"
$arr[1,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$arr[1,11] = 4
$arr[1,12] = 70
$arr[1,13] = 39749
$arr[1,14] = "X"
$arr[1,15] = 39749
$arr[1,16] = "2025-04-15 00:00:00"
$arr[1,17] = -0.3623372395833334
$arr[2,0] = 251391
$arr[2,1] = "BIMEC 2"
$arr[2,2] = 17
$arr[2,3] = 91.640625
$arr[2,4] = "2025-05-07 08:41:45"
$arr[2,5] = "2025-05-07 08:58:45"
$arr[2,6] = "2025-05-07 08:58:45"
$arr[2,7] = "2025-05-07 10:30:24"
$arr[2,8] = 5865
$arr[2,9] = "bobina"
$arr[2,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$arr[2,11] = 5
$arr[2,12] = 70
$arr[2,13] = 39749
$arr[2,14] = "X"
$arr[2,15] = 39749
$arr[2,16] = "2025-04-23 00:00:00"
$arr[2,17] = -0.4377821180555556
$arr[3,0] = 251395
$arr[3,1] = "BIMEC 2"
$arr[3,2] = 17
$arr[3,3] = 35.34375
$arr[3,4] = "2025-05-07 10:30:24"
$arr[3,5] = "2025-05-07 10:47:24"
$arr[3,6] = "2025-05-07 10:47:24"
$arr[3,7] = "2025-05-07 11:22:45"
$arr[3,8] = 2262
$arr[3,9] = "bobina"
$arr[3,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$arr[3,11] = 6
$arr[3,12] = 70
$arr[3,13] = 39749
$arr[3,14] = "X"
$arr[3,15] = 39749
$arr[3,16] = "2025-04-23 00:00:00"
$arr[3,17] = -0.4741319444444445
$arr[4,0] = 251371
$arr[4,1] = "BIMEC 2"
$arr[4,2] = 19
$arr[4,3] = 0
$arr[4,4] = "2025-05-07 11:22:45"
$arr[4,5] = "2025-05-07 11:41:45"
$arr[4,6] = "2025-05-07 11:41:45"
$arr[4,7] = "2025-05-07 11:41:45"
$arr[4,8] = 0
$arr[4,9] = "bobina"
$arr[4,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$arr[4,11] = 4
$arr[4,12] = 70
$arr[4,13] = "39666 (esterno)"
$arr[4,14] = "X"
$arr[4,15] = 39666
$arr[4,16] = "2025-04-24 00:00:00"
$arr[4,17] = -13.48732638888889
$arr[5,0] = 251453
$arr[5,1] = "BIMEC 2"
$arr[5,2] = 17
$arr[5,3] = 78.125
$arr[5,4] = "2025-05-07 11:41:45"
$arr[5,5] = "2025-05-07 11:58:45"
$arr[5,6] = "2025-05-07 11:58:45"
$arr[5,7] = "2025-05-07 13:16:52"
$arr[5,8] = 5000
$arr[5,9] = "bobina"
$arr[5,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$arr[5,11] = 3
$arr[5,12] = 70
$arr[5,13] = "39742 (non in estrazione)"
$arr[5,14] = "X"
$arr[5,15] = 39742
$arr[5,16] = "2025-04-28 00:00:00"
$arr[5,17] = -9.553385416666666
$arr[6,0] = 251396
$arr[6,1] = "BIMEC 2"
$arr[6,2] = 21
$arr[6,3] = 35.34375
$arr[6,4] = "2025-05-07 13:16:52"
$arr[6,5] = "2025-05-07 13:37:52"
$arr[6,6] = "2025-05-07 13:37:52"
$arr[6,7] = "2025-05-07 14:13:13"
$arr[6,8] = 2262
$arr[6,9] = "bobina"
$arr[6,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$arr[6,11] = 6
$arr[6,12] = 70
$arr[6,13] = 39749
$arr[6,14] = "X"
$arr[6,15] = 39749
$arr[6,16] = "2025-05-02 00:00:00"
$arr[6,17] = -0.5925130208333333
$arr[7,0] = 251548
$arr[7,1] = "BIMEC 2"
$arr[7,2] = 19
$arr[7,3] = 206.90625
$arr[7,4] = "2025-05-07 14:13:13"
$arr[7,5] = "2025-05-07 14:32:13"
$arr[7,6] = "2025-05-07 14:32:13"
$arr[7,7] = "2025-05-08 09:59:07"
$arr[7,8] = 13242
$arr[7,9] = "bobina"
$arr[7,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$arr[7,11] = 4
$arr[7,12] = 70
$arr[7,13] = 39749
$arr[7,14] = "X"
$arr[7,15] = 39749
$arr[7,16] = "2025-05-06 00:00:00"
$arr[7,17] = -1.416059027777778
$arr[8,0] = 250923
$arr[8,1] = "BIMEC 2"
$arr[8,2] = 32
$arr[8,3] = 109.46875
$arr[8,4] = "2025-05-08 09:59:07"
$arr[8,5] = "2025-05-08 10:31:07"
$arr[8,6] = "2025-05-08 10:31:07"
$arr[8,7] = "2025-05-08 12:20:35"
$arr[8,8] = 7006
$arr[8,9] = "bobina"
$arr[8,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R9"
$arr[8,11] = 5
$arr[8,12] = 76
$arr[8,13] = 39749
$arr[8,14] = "X"
$arr[8,15] = 39749
$arr[8,16] = "2025-04-07 00:00:00"
$arr[8,17] = -1.514301215277778
$arr[9,0] = 251225
$arr[9,1] = "BIMEC 2"
$arr[9,2] = 17
$arr[9,3] = 0
$arr[9,4] = "2025-05-08 12:20:35"
$arr[9,5] = "2025-05-08 12:37:35"
$arr[9,6] = "2025-05-08 12:37:35"
$arr[9,7] = "2025-05-08 12:37:35"
$arr[9,8] = 0
$arr[9,9] = "bobina"
$arr[9,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R9"
$arr[9,11] = 4
$arr[9,12] = 76
$arr[9,13] = 39747
$arr[9,14] = "X"
$arr[9,15] = 39747
$arr[9,16] = "2025-04-30 00:00:00"
$arr[9,17] = -0.5261067708333333
$arr[10,0] = 251227
$arr[10,1] = "BIMEC 2"
$arr[10,2] = 15
$arr[10,3] = 0
$arr[10,4] = "2025-05-08 12:37:35"
$arr[10,5] = "2025-05-08 12:52:35"
$arr[10,6] = "2025-05-08 12:52:35"
$arr[10,7] = "2025-05-08 12:52:35"
$arr[10,8] = 0
$arr[10,9] = "bobina"
$arr[10,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R9"
$arr[10,11] = 4
$arr[10,12] = 76
$arr[10,13] = 39746
$arr[10,14] = "X"
$arr[10,15] = 39746
$arr[10,16] = "2025-05-05 00:00:00"
$arr[10,17] = -2.5365234375
$arr[11,0] = 251421
$arr[11,1] = "BIMEC 2"
$arr[11,2] = 17
$arr[11,3] = 81.9375
$arr[11,4] = "2025-05-08 12:52:35"
$arr[11,5] = "2025-05-08 13:09:35"
$arr[11,6] = "2025-05-08 13:09:35"
$arr[11,7] = "2025-05-08 14:31:31"
$arr[11,8] = 5244
$arr[11,9] = "bobina"
$arr[11,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R6 ;R9"
$arr[11,11] = 3
$arr[11,12] = 76
$arr[11,13] = "39762 (non in estrazione)"
$arr[11,14] = "X"
$arr[11,15] = 39762
$arr[11,16] = "2025-05-08 00:00:00"
$arr[11,17] = -0.6052300347222223
$arr[12,0] = 251782
$arr[12,1] = "BIMEC 2"
$arr[12,2] = 15
$arr[12,3] = 188.640625
$arr[12,4] = "2025-05-08 14:31:31"
$arr[12,5] = "2025-05-08 14:46:31"
$arr[12,6] = "2025-05-08 14:46:31"
$arr[12,7] = "2025-05-09 09:55:10"
$arr[12,8] = 12073
$arr[12,9] = "bobina"
$arr[12,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R6 ;R9"
$arr[12,11] = 3
$arr[12,12] = 76
$arr[12,13] = 39754
$arr[12,14] = "X"
$arr[12,15] = 39754
$arr[12,16] = "2025-05-16 00:00:00"
$arr[12,17] = -0.4133138020833333
$arr[13,0] = 251050
$arr[13,1] = "R6"
$arr[13,2] = 217
$arr[13,3] = 0
$arr[13,4] = "2025-05-09 07:00:00"
$arr[13,5] = "2025-05-09 10:37:00"
$arr[13,6] = "2025-05-09 10:37:00"
$arr[13,7] = "2025-05-09 10:37:00"
$arr[13,8] = 0
$arr[13,9] = "bobina"
$arr[13,10] = "R6"
$arr[13,11] = 38
$arr[13,12] = 70
$arr[13,13] = 39747
$arr[13,14] = "X"
$arr[13,15] = 39747
$arr[13,16] = "2025-04-16 00:00:00"
$arr[13,17] = -1.442361111111111
$arr[14,0] = 251054
$arr[14,1] = "R6"
$arr[14,2] = 35
$arr[14,3] = 0
$arr[14,4] = "2025-05-09 10:37:00"
$arr[14,5] = "2025-05-09 11:12:00"
$arr[14,6] = "2025-05-09 11:12:00"
$arr[14,7] = "2025-05-09 11:12:00"
$arr[14,8] = 0
$arr[14,9] = "bobina"
$arr[14,10] = "R6"
$arr[14,11] = 38
$arr[14,12] = 70
$arr[14,13] = 39747
$arr[14,14] = "X"
$arr[14,15] = 39747
$arr[14,16] = "2025-04-16 00:00:00"
$arr[14,17] = -1.466666666666667
$arr[15,0] = 251081
$arr[15,1] = "R6"
$arr[15,2] = 125
$arr[15,3] = 42.42253521126761
$arr[15,4] = "2025-05-09 11:12:00"
$arr[15,5] = "2025-05-09 13:17:00"
$arr[15,6] = "2025-05-09 13:17:00"
$arr[15,7] = "2025-05-09 13:59:25"
$arr[15,8] = 3012
$arr[15,9] = "bobina"
$arr[15,10] = "R6"
$arr[15,11] = 20
$arr[15,12] = 70
$arr[15,13] = "39750 (esterno)"
$arr[15,14] = "X"
$arr[15,15] = 39750
$arr[15,16] = "2025-04-23 00:00:00"
$arr[15,17] = -16.58293231612268
$arr[16,0] = 251284
$arr[16,1] = "CASON"
$arr[16,2] = 40.5
$arr[16,3] = 297.0909090909091
$arr[16,4] = "2025-05-09 07:00:00"
$arr[16,5] = "2025-05-09 07:40:30"
$arr[16,6] = "2025-05-09 07:40:30"
$arr[16,7] = "2025-05-09 12:37:35"
$arr[16,8] = 16340
$arr[16,9] = "bobina"
$arr[16,10] = "CASON ;R6"
$arr[16,11] = 7
$arr[16,12] = 70
$arr[16,13] = 39747
$arr[16,14] = "X"
$arr[16,15] = 39747
$arr[16,16] = "2025-05-12 00:00:00"
$arr[16,17] = -1.526104797974537
$arr[17,0] = 251742
$arr[17,1] = "R10"
$arr[17,2] = 30
$arr[17,3] = 134.8524590163935
$arr[17,4] = "2025-05-08 07:00:00"
$arr[17,5] = "2025-05-08 07:30:00"
$arr[17,6] = "2025-05-08 07:30:00"
$arr[17,7] = "2025-05-08 09:44:51"
$arr[17,8] = 8226
$arr[17,9] = "bobina"
$arr[17,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$arr[17,11] = 4
$arr[17,12] = 70
$arr[17,13] = 39749
$arr[17,14] = "X"
$arr[17,15] = 39749
$arr[17,16] = "2025-05-15 00:00:00"
$arr[17,17] = -1.406147540983796
$arr[18,0] = 251840
$arr[18,1] = "R10"
$arr[18,2] = 25
$arr[18,3] = 93.67213114754098
$arr[18,4] = "2025-05-08 09:44:51"
$arr[18,5] = "2025-05-08 10:09:51"
$arr[18,6] = "2025-05-08 10:09:51"
$arr[18,7] = "2025-05-08 11:43:31"
$arr[18,8] = 5714
$arr[18,9] = "bobina"
$arr[18,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$arr[18,11] = 5
$arr[18,12] = 70
$arr[18,13] = 39758
$arr[18,14] = "X"
$arr[18,15] = 39758
$arr[18,16] = "2025-05-09 00:00:00"
$arr[18,17] = -0.4885587431712963
$arr[19,0] = 251456
$arr[19,1] = "R10"
$arr[19,2] = 30
$arr[19,3] = 147.5245901639344
$arr[19,4] = "2025-05-08 11:43:31"
$arr[19,5] = "2025-05-08 12:13:31"
$arr[19,6] = "2025-05-08 12:13:31"
$arr[19,7] = "2025-05-08 14:41:02"
$arr[19,8] = 8999
$arr[19,9] = "bobina"
$arr[19,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$arr[19,11] = 3
$arr[19,12] = 70
$arr[19,13] = 39746
$arr[19,14] = "X"
$arr[19,15] = 39746
$arr[19,16] = "2025-05-09 00:00:00"
$arr[19,17] = -2.611839708564815
$arr[20,0] = 251416
$arr[20,1] = "R10"
$arr[20,2] = 25
$arr[20,3] = 183.9672131147541
$arr[20,4] = "2025-05-08 14:41:02"
$arr[20,5] = "2025-05-09 07:06:02"
$arr[20,6] = "2025-05-09 07:06:02"
$arr[20,7] = "2025-05-09 10:10:00"
$arr[20,8] = 11222
$arr[20,9] = "bobina"
$arr[20,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$arr[20,11] = 2
$arr[20,12] = 70
$arr[20,13] = 39755
$arr[20,14] = 0
$arr[20,15] = 0
$arr[20,16] = "2025-04-23 00:00:00"
$arr[20,17] = 0
$arr[21,0] = 251547
$arr[21,1] = "BIMEC 5"
$arr[21,2] = 34
$arr[21,3] = 184.9154929577465
$arr[21,4] = "2025-05-08 07:00:00"
$arr[21,5] = "2025-05-08 07:34:00"
$arr[21,6] = "2025-05-08 07:34:00"
$arr[21,7] = "2025-05-08 10:38:54"
$arr[21,8] = 13129
$arr[21,9] = "bobina"
$arr[21,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$arr[21,11] = 4
$arr[21,12] = 70
$arr[21,13] = 39749
$arr[21,14] = "X"
$arr[21,15] = 39749
$arr[21,16] = "2025-05-06 00:00:00"
$arr[21,17] = -1.443691314548611
$arr[22,0] = 250759
$arr[22,1] = "BIMEC 5"
$arr[22,2] = 30
$arr[22,3] = 118.2816901408451
$arr[22,4] = "2025-05-08 10:38:54"
$arr[22,5] = "2025-05-08 11:08:54"
$arr[22,6] = "2025-05-08 11:08:54"
$arr[22,7] = "2025-05-08 13:07:11"
$arr[22,8] = 8398
$arr[22,9] = "bobina"
$arr[22,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12"
$arr[22,11] = 4
$arr[22,12] = 76
$arr[22,13] = 39747
$arr[22,14] = "X"
$arr[22,15] = 39747
$arr[22,16] = "2025-03-13 00:00:00"
$arr[22,17] = -0.5466647104861111
$arr[23,0] = 251229
$arr[23,1] = "BIMEC 5"
$arr[23,2] = 34
$arr[23,3] = 263.9295774647887
$arr[23,4] = "2025-05-08 13:07:11"
$arr[23,5] = "2025-05-08 13:41:11"
$arr[23,6] = "2025-05-08 13:41:11"
$arr[23,7] = "2025-05-09 10:05:07"
$arr[23,8] = 18739
$arr[23,9] = "bobina"
$arr[23,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R9"
$arr[23,11] = 6
$arr[23,12] = 70
$arr[23,13] = "39723 (esterno)"
$arr[23,14] = "X"
$arr[23,15] = 39723
$arr[23,16] = "2025-05-15 00:00:00"
$arr[23,17] = 0
$arr[24,0] = 251477
$arr[24,1] = "R12"
$arr[24,2] = 17
$arr[24,3] = 422.5211267605634
$arr[24,4] = "2025-05-08 12:00:00"
$arr[24,5] = "2025-05-08 12:17:00"
$arr[24,6] = "2025-05-08 12:17:00"
$arr[24,7] = "2025-05-09 11:19:31"
$arr[24,8] = 29999
$arr[24,9] = "bobina"
$arr[24,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R6 ;R9"
$arr[24,11] = 3
$arr[24,12] = 76
$arr[24,13] = 39760
$arr[24,14] = "X"
$arr[24,15] = 39760
$arr[24,16] = "2025-04-28 00:00:00"
$arr[24,17] = -2.471889671365741
$arr[25,0] = 251651
$arr[25,1] = "BIMEC 4"
$arr[25,2] = 29
$arr[25,3] = 767.7049180327868
$arr[25,4] = "2025-05-09 07:00:00"
$arr[25,5] = "2025-05-09 07:29:00"
$arr[25,6] = "2025-05-09 07:29:00"
$arr[25,7] = "2025-05-12 12:16:42"
$arr[25,8] = 46830
$arr[25,9] = "bobina"
$arr[25,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R3 ;R6 ;R9"
$arr[25,11] = 5
$arr[25,12] = 76
$arr[25,13] = 39755
$arr[25,14] = 0
$arr[25,15] = 0
$arr[25,16] = "2025-05-12 00:00:00"
$arr[25,17] = 0
$arr[26,0] = 251268
$arr[26,1] = "R3"
$arr[26,2] = 47
$arr[26,3] = 0
$arr[26,4] = "2025-05-08 07:00:00"
$arr[26,5] = "2025-05-08 07:47:00"
$arr[26,6] = "2025-05-08 07:47:00"
$arr[26,7] = "2025-05-08 07:47:00"
$arr[26,8] = 0
$arr[26,9] = "bobina"
$arr[26,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R3 ;R9"
$arr[26,11] = 4
$arr[26,12] = 76
$arr[26,13] = "39666 (non in estrazione)"
$arr[26,14] = "X"
$arr[26,15] = 39666
$arr[26,16] = "2025-04-14 00:00:00"
$arr[26,17] = -24.32430555555555
$arr[27,0] = 251164
$arr[27,1] = "R3"
$arr[27,2] = 47
$arr[27,3] = 204.0816326530612
$arr[27,4] = "2025-05-08 07:47:00"
$arr[27,5] = "2025-05-08 08:34:00"
$arr[27,6] = "2025-05-08 08:34:00"
$arr[27,7] = "2025-05-08 11:58:04"
$arr[27,8] = 10000
$arr[27,9] = "bobina"
$arr[27,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$arr[27,11] = 6
$arr[27,12] = 70
$arr[27,13] = 39749
$arr[27,14] = "X"
$arr[27,15] = 39749
$arr[27,16] = "2025-04-22 00:00:00"
$arr[27,17] = -1.498667800451389
$arr[28,0] = 250894
$arr[28,1] = "R3"
$arr[28,2] = 42
$arr[28,3] = 903.3061224489796
$arr[28,4] = "2025-05-08 11:58:04"
$arr[28,5] = "2025-05-08 12:40:04"
$arr[28,6] = "2025-05-08 12:40:04"
$arr[28,7] = "2025-05-12 11:43:23"
$arr[28,8] = 44262
$arr[28,9] = "bobina"
$arr[28,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R3 ;R6 ;R9"
$arr[28,11] = 5
$arr[28,12] = 76
$arr[28,13] = 39755
$arr[28,14] = 0
$arr[28,15] = 0
$arr[28,16] = "2025-05-05 00:00:00"
$arr[28,17] = 0
$ws.Range("A2:R30").Value = $arr
